$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 21 corresponds to the weather log entry for 2023-11-10 (serial 45240):
# Min Temp, Max Temp, Wind, Inches of Rain updated with corrected readings.
$ws.Range("B21").Value = 43.4
$ws.Range("C21").Value = 51.8
$ws.Range("D21").Value = 13.4
$ws.Range("E21").Value = 0.05
